# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from 45175 (2023-09-06) to 45183 (2023-09-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}
